$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on the cells we are about to write so that
# numeric-looking strings (e.g. "352.70", "2.199.62") are preserved
# exactly as text, matching the source data which stores these as
# inline strings rather than numbers.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '52.151.06'
$ws.Range("E2").Value = '  +0.65%  '
$ws.Range("D3").Value = '2.947.06'
$ws.Range("E3").Value = '  +4.83%  '
$ws.Range("D5").Value = '352.70'
$ws.Range("E5").Value = '  +0.41%  '
$ws.Range("D6").Value = '112.34'
$ws.Range("E6").Value = '  -0.55%  '
$ws.Range("D7").Value = '0.558'
$ws.Range("E7").Value = '  -0.43%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").Value = '0.624'
$ws.Range("E9").Value = '  +0.53%  '
$ws.Range("D10").Value = '39.54'
$ws.Range("E10").Value = '  -1.84%  '
$ws.Range("D11").Value = '0.0880'
$ws.Range("E11").Value = '  +3.50%  '
$ws.Range("E12").Value = '  +0.97%  '
$ws.Range("D13").Value = '20.02'
$ws.Range("E13").Value = '  +0.24%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '3.401.41'
$ws.Range("E14").Value = '  +4.47%  '
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").Value = '7.74'
$ws.Range("E15").Value = '  -0.82%  '
$ws.Range("D16").Value = '2.920.39'
$ws.Range("E16").Value = '  +4.56%  '
$ws.Range("E17").Value = '  +1.30%  '
$ws.Range("D18").Value = '52.213.02'
$ws.Range("E18").Value = '  +0.69%  '
$ws.Range("D19").Value = '7.63'
$ws.Range("E19").Value = '  +0.42%  '
$ws.Range("E20").Value = '  -2.68%  '
$ws.Range("D21").Value = '14.18'
$ws.Range("D22").Value = '0.0₃0980'
$ws.Range("E22").Value = '  +0.49%  '
$ws.Range("D23").Value = '71.09'
$ws.Range("E23").Value = '  +0.63%  '
$ws.Range("D24").Value = '268.29'
$ws.Range("E24").Value = '  -0.15%  '
$ws.Range("E25").Value = '  +0.92%  '
$ws.Range("D26").Value = '0.180'
$ws.Range("E26").Value = '  +10.65%  '
$ws.Range("D27").Value = '27.12'
$ws.Range("E27").Value = '  +3.41%  '
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("D29").Value = '6.97'
$ws.Range("E29").Value = '  +12.62%  '
$ws.Range("D30").Value = '10.62'
$ws.Range("E30").Value = '  +1.12%  '
$ws.Range("D31").Value = '0.103'
$ws.Range("E31").Value = '  +14.30%  '
$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").Value = '36.94'
$ws.Range("E32").Value = '  -4.78%  '
$ws.Range("B33").Value = 'RenderToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D33").Value = '6.03'
$ws.Range("E33").Value = '  +5.41%  '
$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D34").Value = '53.05'
$ws.Range("E34").Value = '  +0.63%  '
$ws.Range("B35").Value = 'Toncoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D35").Value = '2.10'
$ws.Range("E35").Value = '  -7.28%  '
$ws.Range("D36").Value = '0.0453'
$ws.Range("E36").Value = '  -0.23%  '
$ws.Range("D37").Value = '0.998'
$ws.Range("E37").Value = '  -0.15%  '
$ws.Range("D38").Value = '3.39'
$ws.Range("E38").Value = '  +5.91%  '
$ws.Range("D39").Value = '18.65'
$ws.Range("D40").Value = '2.06'
$ws.Range("E40").Value = '  +2.33%  '
$ws.Range("D41").Value = '2.70'
$ws.Range("E41").Value = '  +5.50%  '
$ws.Range("E42").Value = '  +1.20%  '
$ws.Range("D43").Value = '23.27'
$ws.Range("E43").Value = '  +4.35%  '
$ws.Range("E44").Value = '  -2.17%  '
$ws.Range("D45").Value = '2.199.62'
$ws.Range("E45").Value = '  +2.58%  '
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").Value = '2.51'
$ws.Range("E46").Value = '  +0.90%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '3.50'
$ws.Range("E47").Value = '  -0.19%  '
$ws.Range("D48").Value = '110.39'
$ws.Range("E48").Value = '  -9.31%  '
$ws.Range("D49").Value = '0.248'
$ws.Range("E49").Value = '  +9.85%  '
$ws.Range("E50").Value = '  +8.69%  '
$ws.Range("E51").Value = '  -4.38%  '

# Restore the default (Normal) style on the edited range so no stray
# number-format / style index is left behind on these cells.
$ws.Range("B2:E51").Style = "Normal"
